$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, since many new values look numeric
# and Excel would otherwise auto-convert them to numbers, losing the literal
# "x.xxx"-style formatting used throughout this sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.345.93"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.869.20"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "234.83"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "0.4694"
$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").Value = "0.2873"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "0.06573"
$ws.Range("E9").Value = "  +0.39%  "

$ws.Range("D10").Value = "21.64"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").Value = "0.07887"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "96.55"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").Value = "1.865.03"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "0.6926"
$ws.Range("E14").Value = "  +1.54%  "

$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").Value = "268.25"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "30.285.79"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").Value = "14.00"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").Value = "0.000007680"
$ws.Range("E19").Value = "  +3.30%  "

$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "2.110.74"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "5.238"
$ws.Range("E23").Value = "  -1.60%  "

$ws.Range("D24").Value = "6.193"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").Value = "9.403"
$ws.Range("E25").Value = "  +1.89%  "

$ws.Range("D26").Value = "167.54"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").Value = "18.84"
$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D29").Value = "1.361"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("D30").Value = "0.09887"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").Value = "4.396"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("D33").Value = "4.073"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").Value = "0.7026"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "2.724"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("D38").Value = "0.01873"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "2.802"
$ws.Range("E39").Value = "  +6.98%  "

$ws.Range("D40").Value = "6.237"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").Value = "73.41"
$ws.Range("E41").Value = "  -1.72%  "

$ws.Range("D42").Value = "1.953"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4176"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8421"
$ws.Range("E44").Value = "  -0.50%  "

$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").Value = "102.79"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").Value = "967.19"
$ws.Range("E47").Value = "  +1.19%  "

$ws.Range("D48").Value = "7.120"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").Value = "9.106"
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").Value = "34.52"
$ws.Range("E50").Value = "  +1.00%  "

$ws.Range("D51").Value = "0.05678"
